$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NATIONALUM's "Avg. cost" (C3) was updated from 188.49 to 250.5
$ws.Range("C3").Value = 250.5

# The OIL holding (row 4) was removed entirely. Deleting the whole row
# shifts everything below it (the Motilal Oswal Large and Midcap Fund
# row) up by one, so it becomes the new row 4. This also drops the
# "OIL"/"Oil" shared strings since nothing references them anymore.
$ws.Rows.Item(4).Delete()

# Match the final selection recorded in the sheet view
$ws.Range("B4").Select()
